# Insert a new row at row 27 (pushing existing rows 27-99 down to 28-100)
# and populate it with the new weekly record, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(27).Insert()

$ws.Cells.Item(27, 1).Value = 10
$ws.Cells.Item(27, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(27, 3).Value = "La Araucanía"
$ws.Cells.Item(27, 4).Value = 44949
$ws.Cells.Item(27, 5).Value = 9
$ws.Cells.Item(27, 6).Value = 100112030
$ws.Cells.Item(27, 7).Value = "Poroto granado"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 75
$ws.Cells.Item(27, 11).Value = 45000
$ws.Cells.Item(27, 12).Value = 45000
$ws.Cells.Item(27, 13).Value = 45000
$ws.Cells.Item(27, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(27, 15).Value = "Región del Maule"
$ws.Cells.Item(27, 16).Value = 1800
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = "Hortaliza"
